# "Out of PO.xlsx" update — refresh the player / position / team table
# (new players rotated in, positions & teams reassigned).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("De'Andre Hunter",    "SF,PF",    "Cleveland Cavaliers"),
    @("Jaden McDaniels",    "SF,PF",    "Minnesota Timberwolves"),
    @("Bilal Coulibaly",    "SG,SF",    "Washington Wizards"),
    @("Victor Wembanyama",  "C",        "San Antonio Spurs"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Kel'el Ware",        "PF,C",     "Miami Heat"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("Naji Marshall",      "SG,SF",    "Dallas Mavericks"),
    @("Julius Randle",      "PF,C",     "Minnesota Timberwolves"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
    @("Alperen Sengün",     "C",        "Houston Rockets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
